$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.5036099671744694
$ws.Range("C2").Value = 0.0405203490245043
$ws.Range("D2").Value = 0.1878942985744487
$ws.Range("E2").Value = 0.4023142941575912
$ws.Range("F2").Value = 3.48751082828565
$ws.Range("K2").Value = 0.4667827300587533
$ws.Range("B3").Value = 0.4724063787152772
$ws.Range("C3").Value = 0.03569615230428269
$ws.Range("D3").Value = 0.1760989280804068
$ws.Range("E3").Value = 0.3508673324256222
$ws.Range("F3").Value = 3.266122130889642
$ws.Range("K3").Value = 0.4333145590498191
$ws.Range("B4").Value = 0.4537886395940234
$ws.Range("C4").Value = 0.032786681855697
$ws.Range("D4").Value = 0.168813765352354
$ws.Range("E4").Value = 0.3194104324273752
$ws.Range("F4").Value = 3.13049815831107
$ws.Range("K4").Value = 0.4133131964094048
$ws.Range("B5").Value = 0.446336625980166
$ws.Range("C5").Value = 0.03161388729661496
$ws.Range("D5").Value = 0.1658337207018263
$ws.Range("E5").Value = 0.3066218381586054
$ws.Range("F5").Value = 3.075301371088585
$ws.Range("K5").Value = 0.4052987025102368
$ws.Range("B6").Value = 0.4451073337293394
$ws.Range("C6").Value = 0.03141990975703379
$ws.Range("D6").Value = 0.1653381897334327
$ws.Range("E6").Value = 0.3045000542300045
$ws.Range("F6").Value = 3.066140083660343
$ws.Range("K6").Value = 0.4039760823098106
$ws.Range("B7").Value = 0.4536875947337649
$ws.Range("C7").Value = 0.0327708136670708
$ws.Range("D7").Value = 0.1687736217987492
$ws.Range("E7").Value = 0.3192378417100628
$ws.Range("F7").Value = 3.129753476227449
$ws.Range("K7").Value = 0.4132045607218515
$ws.Range("B8").Value = 0.4927378209911524
$ws.Range("C8").Value = 0.03884580997198839
$ws.Range("D8").Value = 0.1838358066363242
$ws.Range("E8").Value = 0.3845463406986056
$ws.Range("F8").Value = 3.411107596284126
$ws.Range("K8").Value = 0.4551281166546914
$ws.Range("B9").Value = 0.5736739766611549
$ws.Range("C9").Value = 0.05119463413410585
$ws.Range("D9").Value = 0.2130585166598422
$ws.Range("E9").Value = 0.513793844763029
$ws.Range("F9").Value = 3.965631833296101
$ws.Range("K9").Value = 0.5417697655767029
$ws.Range("B10").Value = 0.6358907497865118
$ws.Range("C10").Value = 0.06055955985119965
$ws.Range("D10").Value = 0.2343733135178923
$ws.Range("E10").Value = 0.6096635881758914
$ws.Range("F10").Value = 4.375243652213783
$ws.Range("K10").Value = 0.6082450606347436
$ws.Range("B11").Value = 0.6648141964094521
$ws.Range("C11").Value = 0.06488927881406426
$ws.Range("D11").Value = 0.244044286443625
$ws.Range("E11").Value = 0.6535182957296115
$ws.Range("F11").Value = 4.562177163855381
$ws.Range("K11").Value = 0.6391254755381794
$ws.Range("B12").Value = 0.675857676010736
$ws.Range("C12").Value = 0.06653930966280086
$ws.Range("D12").Value = 0.2477034236533484
$ws.Range("E12").Value = 0.6701632973674663
$ws.Range("F12").Value = 4.633058606358702
$ws.Range("K12").Value = 0.65091327339087
$ws.Range("B13").Value = 0.6734752019108612
$ws.Range("C13").Value = 0.0661834734066673
$ws.Range("D13").Value = 0.2469154893751977
$ws.Range("E13").Value = 0.6665767486054506
$ws.Range("F13").Value = 4.617788725657761
$ws.Range("K13").Value = 0.6483703454628369
$ws.Range("B14").Value = 0.6657209214455975
$ws.Range("C14").Value = 0.06502481555871498
$ws.Range("D14").Value = 0.244345383662079
$ws.Range("E14").Value = 0.6548869062502405
$ws.Range("F14").Value = 4.568006699117689
$ws.Range("K14").Value = 0.6400933683955827
$ws.Range("B15").Value = 0.6609830716624856
$ws.Range("C15").Value = 0.06431648011036373
$ws.Range("D15").Value = 0.2427707397322934
$ws.Range("E15").Value = 0.6477316085806422
$ws.Range("F15").Value = 4.537526212058935
$ws.Range("K15").Value = 0.6350357885033304
$ws.Range("B16").Value = 0.6340131557641371
$ws.Range("C16").Value = 0.06027804122604152
$ws.Range("D16").Value = 0.2337408342029619
$ws.Range("E16").Value = 0.6068027451735958
$ws.Range("F16").Value = 4.3630398117634
$ws.Range("K16").Value = 0.6062400069797604
$ws.Range("B17").Value = 0.617628051291689
$ws.Range("C17").Value = 0.05781872960029943
$ws.Range("D17").Value = 0.2281952259481557
$ws.Range("E17").Value = 0.5817587434397922
$ws.Range("F17").Value = 4.256156879186364
$ws.Range("K17").Value = 0.5887401450713696
$ws.Range("B18").Value = 0.6082621038126774
$ws.Range("C18").Value = 0.05641071874919135
$ws.Range("D18").Value = 0.2250031285520748
$ws.Range("E18").Value = 0.5673767794070557
$ws.Range("F18").Value = 4.194736563425209
$ws.Range("K18").Value = 0.5787348517135342
$ws.Range("B19").Value = 0.6051009270975385
$ws.Range("C19").Value = 0.0559350965473584
$ws.Range("D19").Value = 0.2239219077475951
$ws.Range("E19").Value = 0.5625110923185304
$ws.Range("F19").Value = 4.173950073580272
$ws.Range("K19").Value = 0.5753575057075579
$ws.Range("B20").Value = 0.6193662255468269
$ws.Range("C20").Value = 0.05807984943244549
$ws.Range("D20").Value = 0.2287858106901695
$ws.Range("E20").Value = 0.5844223494058411
$ws.Range("F20").Value = 4.267528907387884
$ws.Range("K20").Value = 0.5905967944376584
$ws.Range("B21").Value = 0.6679960649001657
$ws.Range("C21").Value = 0.06536485388205904
$ws.Range("D21").Value = 0.2451003641551779
$ws.Range("E21").Value = 0.658319434249762
$ws.Range("F21").Value = 4.582626280412512
$ws.Range("K21").Value = 0.6425219487592813
$ws.Range("B22").Value = 0.7003084580388474
$ws.Range("C22").Value = 0.07018721421121654
$ws.Range("D22").Value = 0.2557453611652534
$ws.Range("E22").Value = 0.7068395096332125
$ws.Range("F22").Value = 4.789110879957605
$ws.Range("K22").Value = 0.6770072031174266
$ws.Range("B23").Value = 0.6830137140461545
$ws.Range("C23").Value = 0.06760767771636722
$ws.Range("D23").Value = 0.2500653350336108
$ws.Range("E23").Value = 0.6809218496354106
$ws.Range("F23").Value = 4.678853257815604
$ws.Range("K23").Value = 0.6585508671337834
$ws.Range("B24").Value = 0.6185802286742899
$ws.Range("C24").Value = 0.05796177887540921
$ws.Range("D24").Value = 0.2285188193743579
$ws.Range("E24").Value = 0.5832180830372948
$ws.Range("F24").Value = 4.262387526469752
$ws.Range("K24").Value = 0.5897572301927596
$ws.Range("B25").Value = 0.5513021923536314
$ws.Range("C25").Value = 0.04780457923651227
$ws.Range("D25").Value = 0.2051828911296525
$ws.Range("E25").Value = 0.4786830788613656
$ws.Range("F25").Value = 3.815268189410091
$ws.Range("K25").Value = 0.5178445011081578
